$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) and "全部类型" (sheet4) share the same row layout,
# except sheet4 has one extra row (row 9) inserted for a concert entry that
# does not exist on "展览". Apply the matching edits to both sheets.

$sheetExhibit = $wb.Worksheets.Item("展览")
$sheetAll = $wb.Worksheets.Item("全部类型")

# G2: 最低票价 22.33 -> "不可售" (ticket no longer sellable)
$sheetExhibit.Range("G2").Value = "不可售"
$sheetAll.Range("G2").Value = "不可售"

# F column (想去人数 / want-to-go count) updates on "展览" sheet
$sheetExhibit.Range("F6").Value = 35
$sheetExhibit.Range("F7").Value = 568
$sheetExhibit.Range("F9").Value = 51
$sheetExhibit.Range("F12").Value = 2044
$sheetExhibit.Range("F14").Value = 142
$sheetExhibit.Range("F15").Value = 1328
$sheetExhibit.Range("F16").Value = 464
$sheetExhibit.Range("F17").Value = 20
$sheetExhibit.Range("F23").Value = 47
$sheetExhibit.Range("F26").Value = 1100
$sheetExhibit.Range("F28").Value = 330
$sheetExhibit.Range("F29").Value = 176
$sheetExhibit.Range("F30").Value = 267
$sheetExhibit.Range("F31").Value = 318

# F column updates on "全部类型" sheet (rows offset by +1 from row 9 on)
$sheetAll.Range("F6").Value = 35
$sheetAll.Range("F7").Value = 568
$sheetAll.Range("F10").Value = 51
$sheetAll.Range("F13").Value = 2044
$sheetAll.Range("F15").Value = 142
$sheetAll.Range("F16").Value = 1328
$sheetAll.Range("F17").Value = 464
$sheetAll.Range("F18").Value = 20
$sheetAll.Range("F24").Value = 47
$sheetAll.Range("F27").Value = 1100
$sheetAll.Range("F29").Value = 330
$sheetAll.Range("F30").Value = 176
$sheetAll.Range("F31").Value = 267
$sheetAll.Range("F32").Value = 318
